$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 49 - copy style/formatting of row 48's A cell (bold/bordered/centered)
# so the new index cell A49 matches the look of A2:A48.
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A49").Value = 47

$data = New-Object 'object[,]' 48,3
$data[0,0] = 0.1966745697613365
$data[0,1] = 3.398016581806534
$data[0,2] = 0.3479857034036495
$data[1,0] = 4.288805417111569
$data[1,1] = 7.176095580056504
$data[1,2] = 0.6570317406796442
$data[2,0] = 7.158434484609617
$data[2,1] = 6.054386099514097
$data[2,2] = 0.3595808969825859
$data[3,0] = 10.68787202622885
$data[3,1] = 2.471546475542326
$data[3,2] = 0.4931104655600663
$data[4,0] = 11.13035619344921
$data[4,1] = 5.416843953226136
$data[4,2] = 0.6790733474543994
$data[5,0] = 11.23509330811767
$data[5,1] = 7.643562011345336
$data[5,2] = 0.4735205047252736
$data[6,0] = 14.23300387429347
$data[6,1] = 4.000594141629424
$data[6,2] = 0.6564052726327982
$data[7,0] = 16.57644194693058
$data[7,1] = 7.820603897435421
$data[7,2] = 0.5493231789877323
$data[8,0] = 16.85118715735318
$data[8,1] = 2.208994526369969
$data[8,2] = 0.2519192563578154
$data[9,0] = 18.75890062103131
$data[9,1] = 3.452063717791699
$data[9,2] = 0.5305743818499562
$data[10,0] = 19.73781561500588
$data[10,1] = 2.972830468257571
$data[10,2] = 0.1930303733771813
$data[11,0] = 21.53756999053554
$data[11,1] = 2.762419706364087
$data[11,2] = 0.3161397342912934
$data[12,0] = 24.21904830946389
$data[12,1] = 5.977255868190415
$data[12,2] = 0.5135663869071035
$data[13,0] = 26.33950580858796
$data[13,1] = 2.269286081930728
$data[13,2] = 0.4175729612007791
$data[14,0] = 30.11432172914718
$data[14,1] = 3.502639979772928
$data[14,2] = 0.2885986895624441
$data[15,0] = 38.01357964048826
$data[15,1] = 5.687104863253512
$data[15,2] = 0.8249750860183273
$data[16,0] = 39.40097672098837
$data[16,1] = 2.840729492429538
$data[16,2] = 0.1526213693334381
$data[17,0] = 40.6309559186617
$data[17,1] = 2.648039142228023
$data[17,2] = 0.4995017264992364
$data[18,0] = 42.01788358052395
$data[18,1] = 9.032533005398113
$data[18,2] = 0.4584354274854955
$data[19,0] = 42.47961878584243
$data[19,1] = 3.759912426431518
$data[19,2] = 0.5526496099353536
$data[20,0] = 48.21154182148152
$data[20,1] = 4.884695167223111
$data[20,2] = 0.6551279650008076
$data[21,0] = 51.87343159317665
$data[21,1] = 5.1016678442904
$data[21,2] = 0.5111985910070951
$data[22,0] = 52.6826084417166
$data[22,1] = 4.734940375280638
$data[22,2] = 0.4836290784663029
$data[23,0] = 53.2225778079198
$data[23,1] = 3.237682982728692
$data[23,2] = 0.2925664639357802
$data[24,0] = 53.37499097221539
$data[24,1] = 4.113022030632496
$data[24,2] = 0.6839281535466319
$data[25,0] = 53.44336866090006
$data[25,1] = 4.917759825717884
$data[25,2] = 0.8125443301245073
$data[26,0] = 55.16454869715721
$data[26,1] = 3.988384839021643
$data[26,2] = 0.8016877018594017
$data[27,0] = 55.92935279552461
$data[27,1] = 2.434064770685969
$data[27,2] = 0.298923147838203
$data[28,0] = 56.96420929969987
$data[28,1] = 4.504094893883845
$data[28,2] = 0.2375575181956196
$data[29,0] = 57.79034672805739
$data[29,1] = 5.837571719407417
$data[29,2] = 0.7354759494416794
$data[30,0] = 61.10037541310411
$data[30,1] = 4.875486909724401
$data[30,2] = 0.5919267840021381
$data[31,0] = 62.89973489271011
$data[31,1] = 3.168613224775923
$data[31,2] = 0.2023426971452527
$data[32,0] = 67.03346764765708
$data[32,1] = 9.145873471740076
$data[32,2] = 0.6045041708393477
$data[33,0] = 69.12075617778238
$data[33,1] = 3.867017836362711
$data[33,2] = 0.4997329595090217
$data[34,0] = 72.21893476295347
$data[34,1] = 6.335309809597432
$data[34,2] = 0.7107818045511816
$data[35,0] = 73.36888922488927
$data[35,1] = 4.280088842762563
$data[35,2] = 0.3860151713919697
$data[36,0] = 75.18176370430578
$data[36,1] = 7.727121116208082
$data[36,2] = 0.3711330171274874
$data[37,0] = 76.25327634573185
$data[37,1] = 4.063568127503721
$data[37,2] = 0.7214592194034756
$data[38,0] = 82.65041306726405
$data[38,1] = 2.177300345075906
$data[38,2] = 0.1526963793474064
$data[39,0] = 84.25923031779921
$data[39,1] = 6.756571533362284
$data[39,2] = 0.7380599081193271
$data[40,0] = 84.52736561070161
$data[40,1] = 2.827621651999604
$data[40,2] = 0.5407362699526078
$data[41,0] = 85.69280114665088
$data[41,1] = 5.399941233491229
$data[41,2] = 0.6774089583142527
$data[42,0] = 88.36869742799929
$data[42,1] = 6.428142611106714
$data[42,2] = 0.9045113553892575
$data[43,0] = 90.05438927264224
$data[43,1] = 7.0372951426153
$data[43,2] = 0.7317701132050558
$data[44,0] = 96.00512675544741
$data[44,1] = 7.455404693276503
$data[44,2] = 0.682022596229854
$data[45,0] = 97.37350423462532
$data[45,1] = 4.691439591554865
$data[45,2] = 0.4997150791052916
$data[46,0] = 98.0518589920964
$data[46,1] = 3.156895405658013
$data[46,2] = 0.3349985056259042
$data[47,0] = 99.99211251130967
$data[47,1] = 2.503179914633444
$data[47,2] = 0.2052234198103534

$ws.Range("B2:D49").Value = $data
